$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates to the cryptos price/volume table (scraped refresh).
# D-column price cells are forced to Text via a leading quote-prefix so that
# numeric-looking strings (e.g. "25.00", "243.11") keep their exact source
# formatting instead of being auto-coerced to Number by Excel.

$ws.Range("D2").Value = '''29.377.10'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '''1.880.66'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '''0.7224'
$ws.Range("E5").Value = '  +1.70%  '

$ws.Range("D6").Value = '''243.11'

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").Value = '''0.08029'
$ws.Range("E8").Value = '  +2.88%  '

$ws.Range("D9").Value = '''0.3146'
$ws.Range("E9").Value = '  +1.17%  '

$ws.Range("D10").Value = '''25.00'
$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("D11").Value = '''0.08170'
$ws.Range("E11").Value = '  -3.04%  '

$ws.Range("D12").Value = '''1.916.98'
$ws.Range("E12").Value = '  +2.83%  '

$ws.Range("D13").Value = '''94.66'
$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").Value = '''5.234'
$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("D15").Value = '''0.7124'
$ws.Range("E15").Value = '  -0.16%  '

$ws.Range("D16").Value = '''6.427'
$ws.Range("E16").Value = '  +5.59%  '

$ws.Range("D17").Value = '''0.000008484'
$ws.Range("E17").Value = '  +2.00%  '

$ws.Range("D18").Value = '''29.375.77'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").Value = '''244.52'
$ws.Range("E19").Value = '  +1.61%  '

$ws.Range("D20").Value = '''13.33'
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").Value = '''2.127.01'
$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").Value = '''1.002'
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("D23").Value = '''7.767'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("E24").Value = '  +0.14%  '

$ws.Range("D25").Value = '''0.1607'
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("D26").Value = '''162.70'
$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("D27").Value = '''9.047'
$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '''18.54'
$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").Value = '''4.405'
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Value = '''4.283'
$ws.Range("E31").Value = '  -1.01%  '

$ws.Range("D32").Value = '''1.241'
$ws.Range("E32").Value = '  -4.58%  '

$ws.Range("D33").Value = '''0.05356'
$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("D34").Value = '''1.941'
$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("D35").Value = '''0.7637'
$ws.Range("E35").Value = '  +1.71%  '

$ws.Range("D36").Value = '''1.178'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").Value = '''2.701'
$ws.Range("E37").Value = '  +0.33%  '

$ws.Range("D38").Value = '''0.01873'
$ws.Range("E38").Value = '  -0.51%  '

$ws.Range("D39").Value = '''1.261.51'
$ws.Range("E39").Value = '  +2.87%  '

$ws.Range("D40").Value = '''2.766'
$ws.Range("E40").Value = '  +1.42%  '

$ws.Range("D41").Value = '''6.444'
$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("D42").Value = '''113.60'
$ws.Range("E42").Value = '  +3.97%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''74.47'
$ws.Range("E43").Value = '  +2.80%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.9059'
$ws.Range("E44").Value = '  +1.25%  '

$ws.Range("E45").Value = '  +6.57%  '

$ws.Range("E46").Value = '  +0.25%  '

$ws.Range("D47").Value = '''2.029.27'
$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("D48").Value = '''1.803'
$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").Value = '''0.5199'
$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("D50").Value = '''9.520'
$ws.Range("E50").Value = '  +0.85%  '

$ws.Range("D51").Value = '''0.4340'
$ws.Range("E51").Value = '  +0.33%  '
